# Update the cryptocurrency price/volume symbol list.
# These D-column "Price" cells are stored as text (inlineStr) in the workbook,
# so we force Text number-formatting before assigning the new values to stop
# Excel from auto-converting the numeric-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    "D2"  = "244.05"
    "D3"  = "25.01"
    "D4"  = "5.172"
    "D5"  = "0.05686"
    "D6"  = "6.515"
    "D7"  = "2.958"
    "D8"  = "0.8090"
    "D9"  = "0.8394"
    "D10" = "0.1331"
    "D11" = "0.06959"
    "D12" = "0.02828"
    "D14" = "0.001519"
    "D15" = "0.0005999"
    "D16" = "0.006161"
    "D17" = "3.499"
    "D18" = "2.092"
    "D20" = "0.03207"
    "D22" = "3.737"
    "D23" = "0.04673"
    "D25" = "0.001236"
    "D26" = "0.004241"
    "D27" = "0.00009691"
    "D28" = "0.0001966"
    "D40" = "0.03623"
    "D41" = "0.006284"
    "D42" = "0.1050"
    "D43" = "0.002710"
    "D44" = "0.007376"
    "D45" = "0.00005299"
    "D48" = "0.002292"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) column (E) text updates
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
